# Abstinence.xlsx — append 20 new daily-tracking rows (182-201) to the
# "records" table, extending it from A4:N181 to A4:N201.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row data: date serial (E) and the five tracked metrics F..J
$newData = @(
    @(182, 44061, 0, 0, 0, 0, 0),
    @(183, 44062, 0, 0, 0, 0, 0),
    @(184, 44063, 0, 0, 1, 1, 0),
    @(185, 44064, 0, 0, 0, 0, 0),
    @(186, 44065, 0, 0, 0, 0, 0),
    @(187, 44066, 0, 0, 0, 0, 0),
    @(188, 44067, 0, 0, 0, 0, 0),
    @(189, 44068, 0, 0, 0, 0, 0),
    @(190, 44069, 0, 0, 0, 0, 0),
    @(191, 44070, 0, 0, 0, 0, 0),
    @(192, 44071, 0, 0, 0, 0, 0),
    @(193, 44072, 0, 0, 0, 0, 0),
    @(194, 44073, 0, 0, 0, 0, 0),
    @(195, 44074, 0, 0, 0, 0, 0),
    @(196, 44075, 0, 0, 0, 0, 0),
    @(197, 44076, 0, 0, 0, 0, 0),
    @(198, 44077, 0, 0, 0, 0, 0),
    @(199, 44078, 0, 0, 0, 0, 0),
    @(200, 44079, 0, 0, 0, 0, 0),
    @(201, 44080, 0, 0, 0, 0, 0)
)

foreach ($rowData in $newData) {
    $r = $rowData[0]
    $prev = $r - 1

    # E..J: date + the five raw metric columns (plain values)
    $ws.Cells.Item($r, 5).Value = $rowData[1]
    $ws.Cells.Item($r, 6).Value = $rowData[2]
    $ws.Cells.Item($r, 7).Value = $rowData[3]
    $ws.Cells.Item($r, 8).Value = $rowData[4]
    $ws.Cells.Item($r, 9).Value = $rowData[5]
    $ws.Cells.Item($r, 10).Value = $rowData[6]

    # A..D: calculated-column formulas (same formulas used by the table)
    $ws.Cells.Item($r, 1).Formula = '=YEAR(records[[#This Row],[日期]])&" 年"'
    $ws.Cells.Item($r, 2).Formula = '="第 "&INT(MONTH(records[[#This Row],[日期]])/3+1)&" 季度"'
    $ws.Cells.Item($r, 3).Formula = '=MONTH(records[[#This Row],[日期]])&" 月"'
    $ws.Cells.Item($r, 4).Formula = '="第 "&WEEKNUM(records[[#This Row],[日期]],2)&" 周"'

    # L..M: running streak counters, each referencing the row above
    $ws.Cells.Item($r, 12).Formula = "=IF(((records[Porn-Video]+records[Masturbation]+records[Sexual-Intercourse])>0), 0, L$prev+1)"
    $ws.Cells.Item($r, 13).Formula = "=IF(((records[Porn-Video]+records[Masturbation]+records[Sexual-Intercourse])>0), M$prev+1, 0)"
}

# Grow the table ("records") so the new rows become part of it (keeps
# headerRow/autoFilter/calculated columns consistent).
$lo = $ws.ListObjects.Item(1)
$newTableRange = $ws.Range("A4:N201")
$lo.Resize($newTableRange)

# Match the saved view state from the diff (scrolled down, new selection).
$ws.Application.ActiveWindow.ScrollRow = 166
$sheetView = $ws.Range("A166")
$ws.Cells.Item(199, 9).Select()
